# Scheduled runner update: refresh profit-calc figures on the per-job
# "Atomos_Profits" sheets (current average price / leve profit columns).
# Columns: H=currentAveragePrice, I=currentAveragePriceNQ, J=currentAveragePriceHQ,
#          K=LevePriceNQ, L=LevePriceHQ, M=LeveProfitNQ, N=LeveProfitHQ

$wb = $excel.ActiveWorkbook

# ---- ALC ----
$ws = $wb.Worksheets.Item("ALC")

$ws.Range("H17").Value = 1255.6964
$ws.Range("J17").Value = 1255.6964
$ws.Range("L17").Value = 3767.0892
$ws.Range("N17").Value = -4103.0892

$ws.Range("H20").Value = 35512
$ws.Range("I20").Value = 1000
$ws.Range("J20").Value = 70024
$ws.Range("K20").Value = 1000
$ws.Range("L20").Value = 70024
$ws.Range("M20").Value = -770
$ws.Range("N20").Value = -70484

$ws.Range("H33").Value = 567.88464
$ws.Range("I33").Value = 417.625
$ws.Range("K33").Value = 417.625
$ws.Range("M33").Value = -188.625

$ws.Range("H34").Value = 83335410
$ws.Range("I34").Value = 111111870
$ws.Range("J34").Value = 6016.3335
$ws.Range("K34").Value = 111111870
$ws.Range("L34").Value = 6016.3335
$ws.Range("M34").Value = -111111667
$ws.Range("N34").Value = -6422.3335

$ws.Range("H35").Value = 35512
$ws.Range("I35").Value = 1000
$ws.Range("J35").Value = 70024
$ws.Range("K35").Value = 1000
$ws.Range("L35").Value = 70024
$ws.Range("M35").Value = -621
$ws.Range("N35").Value = -70782

$ws.Range("H36").Value = 83335410
$ws.Range("I36").Value = 111111870
$ws.Range("J36").Value = 6016.3335
$ws.Range("K36").Value = 111111870
$ws.Range("L36").Value = 6016.3335
$ws.Range("M36").Value = -111111155
$ws.Range("N36").Value = -7446.3335

$ws.Range("H62").Value = 3195.9375
$ws.Range("I62").Value = 2398.2222
$ws.Range("J62").Value = 4221.5713
$ws.Range("K62").Value = 2398.2222
$ws.Range("L62").Value = 4221.5713
$ws.Range("M62").Value = -1774.2222
$ws.Range("N62").Value = -5469.5713

$ws.Range("H65").Value = 3195.9375
$ws.Range("I65").Value = 2398.2222
$ws.Range("J65").Value = 4221.5713
$ws.Range("K65").Value = 11991.111
$ws.Range("L65").Value = 21107.8565
$ws.Range("M65").Value = -8871.111000000001
$ws.Range("N65").Value = -27347.8565

$ws.Range("H98").Value = 3904.2354
$ws.Range("I98").Value = 2721
$ws.Range("J98").Value = 7749.75
$ws.Range("K98").Value = 2721
$ws.Range("L98").Value = 7749.75
$ws.Range("M98").Value = -1223
$ws.Range("N98").Value = -10745.75

$ws.Range("H122").Value = 3904.2354
$ws.Range("I122").Value = 2721
$ws.Range("J122").Value = 7749.75
$ws.Range("K122").Value = 8163
$ws.Range("L122").Value = 23249.25
$ws.Range("M122").Value = -5713
$ws.Range("N122").Value = -28149.25

# ---- ARM ----
$ws = $wb.Worksheets.Item("ARM")

$ws.Range("H39").Value = 9008
$ws.Range("I39").Value = 9008
$ws.Range("J39").Value = 0
$ws.Range("K39").Value = 9008
$ws.Range("L39").Value = 0
$ws.Range("M39").Value = -8488
$ws.Range("N39").ClearContents()

$ws.Range("H63").Value = 3668.5
$ws.Range("I63").Value = 2335.625
$ws.Range("J63").Value = 9000
$ws.Range("K63").Value = 2335.625
$ws.Range("L63").Value = 9000
$ws.Range("M63").Value = -1649.625
$ws.Range("N63").Value = -10372

$ws.Range("H66").Value = 3668.5
$ws.Range("I66").Value = 2335.625
$ws.Range("J66").Value = 9000
$ws.Range("K66").Value = 11678.125
$ws.Range("L66").Value = 45000
$ws.Range("M66").Value = -8246.125
$ws.Range("N66").Value = -51864

$ws.Range("H74").Value = 1202.5
$ws.Range("I74").Value = 936.6667
$ws.Range("J74").Value = 2000
$ws.Range("K74").Value = 936.6667
$ws.Range("L74").Value = 2000
$ws.Range("M74").Value = -62.66669999999999
$ws.Range("N74").Value = -3748

$ws.Range("H77").Value = 1202.5
$ws.Range("I77").Value = 936.6667
$ws.Range("J77").Value = 2000
$ws.Range("K77").Value = 4683.3335
$ws.Range("L77").Value = 10000
$ws.Range("M77").Value = -315.3334999999997
$ws.Range("N77").Value = -18736

$ws.Range("H134").Value = 71291.60000000001
$ws.Range("J134").Value = 71291.60000000001
$ws.Range("L134").Value = 71291.60000000001
$ws.Range("N134").Value = -81431.60000000001

# ---- BSM ----
$ws = $wb.Worksheets.Item("BSM")

$ws.Range("H38").Value = 70036
$ws.Range("J38").Value = 70036
$ws.Range("L38").Value = 70036
$ws.Range("N38").Value = -70868

$ws.Range("H134").Value = 3408.4482
$ws.Range("I134").Value = 3241.087
$ws.Range("J134").Value = 4050
$ws.Range("K134").Value = 9723.261
$ws.Range("L134").Value = 12150
$ws.Range("M134").Value = -7188.261
$ws.Range("N134").Value = -17220

# ---- CUL ----
$ws = $wb.Worksheets.Item("CUL")

$ws.Range("H3").Value = 4191.8
$ws.Range("I3").Value = 1489.75
$ws.Range("J3").Value = 15000
$ws.Range("K3").Value = 4469.25
$ws.Range("L3").Value = 45000
$ws.Range("M3").Value = -4357.25
$ws.Range("N3").Value = -45224

$ws.Range("H18").Value = 562.63635
$ws.Range("I18").Value = 229.875
$ws.Range("J18").Value = 1450
$ws.Range("K18").Value = 689.625
$ws.Range("L18").Value = 4350
$ws.Range("M18").Value = -520.625
$ws.Range("N18").Value = -4688

$ws.Range("H134").Value = 2428.4119
$ws.Range("I134").Value = 1069
$ws.Range("J134").Value = 3957.75
$ws.Range("K134").Value = 3207
$ws.Range("L134").Value = 11873.25
$ws.Range("M134").Value = 1863
$ws.Range("N134").Value = -22013.25

$ws.Range("H140").Value = 12852728
$ws.Range("I140").Value = 15188543
$ws.Range("J140").Value = 5740
$ws.Range("K140").Value = 45565629
$ws.Range("L140").Value = 17220
$ws.Range("M140").Value = -45560449
$ws.Range("N140").Value = -27580

# ---- LTW ----
$ws = $wb.Worksheets.Item("LTW")

$ws.Range("H61").Value = 76925110
$ws.Range("I61").Value = 83334700
$ws.Range("J61").Value = 10000
$ws.Range("K61").Value = 83334700
$ws.Range("L61").Value = 10000
$ws.Range("M61").Value = -83334498
$ws.Range("N61").Value = -10404

$ws.Range("H113").Value = 76925110
$ws.Range("I113").Value = 83334700
$ws.Range("J113").Value = 10000
$ws.Range("K113").Value = 83334700
$ws.Range("L113").Value = 10000
$ws.Range("M113").Value = -83332530
$ws.Range("N113").Value = -14340

# ---- WVR ----
$ws = $wb.Worksheets.Item("WVR")

$ws.Range("H81").Value = 1167.6666
$ws.Range("I81").Value = 1000.5
$ws.Range("J81").Value = 1502
$ws.Range("K81").Value = 2001
$ws.Range("L81").Value = 3004
$ws.Range("M81").Value = -940
$ws.Range("N81").Value = -5126

$ws.Range("H84").Value = 1167.6666
$ws.Range("I84").Value = 1000.5
$ws.Range("J84").Value = 1502
$ws.Range("K84").Value = 10005
$ws.Range("L84").Value = 15020
$ws.Range("M84").Value = -4701
$ws.Range("N84").Value = -25628

$ws.Range("H96").Value = 21092.4
$ws.Range("I96").Value = 1365.5
$ws.Range("K96").Value = 1365.5
$ws.Range("M96").Value = 7.5

$ws.Range("H122").Value = 287383.72
$ws.Range("I122").Value = 346291.97
$ws.Range("J122").Value = 2660.5
$ws.Range("K122").Value = 1038875.91
$ws.Range("L122").Value = 7981.5
$ws.Range("M122").Value = -1036425.91
$ws.Range("N122").Value = -12881.5
